# Populate newly scraped article URLs into the "AI" and "遠隔会議" sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1 = "AI"
$wsAI = $wb.Worksheets.Item(1)
$aiUrls = @(
    "https://www.itmedia.co.jp/news/articles/1912/06/news103.html",
    "https://www.itmedia.co.jp/news/articles/1912/06/news111.html",
    "https://www.itmedia.co.jp/business/articles/1912/06/news133.html",
    "https://www.itmedia.co.jp/news/articles/1912/06/news036.html",
    "https://www.itmedia.co.jp/business/articles/1912/06/news022.html",
    "https://www.itmedia.co.jp/news/articles/1912/06/news060.html",
    "https://www.itmedia.co.jp/news/articles/1912/06/news086.html",
    "https://www.itmedia.co.jp/business/articles/1912/05/news017.html",
    "https://www.itmedia.co.jp/business/articles/1912/05/news043.html",
    "https://www.itmedia.co.jp/news/articles/1912/05/news119.html"
)

$row = 3
foreach ($url in $aiUrls) {
    $cell = $wsAI.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row = $row + 1
}

# Sheet 8 = "遠隔会議" (remote conference)
$wsRemote = $wb.Worksheets.Item(8)
$remoteUrls = @(
    "https://www.itmedia.co.jp/news/articles/1911/08/news144.html",
    "https://www.itmedia.co.jp/business/articles/1911/01/news005.html",
    "https://www.itmedia.co.jp/news/articles/1910/03/news136.html",
    "https://www.itmedia.co.jp/news/articles/1902/26/news005.html",
    "https://www.itmedia.co.jp/news/articles/1812/05/news005.html",
    "https://www.itmedia.co.jp/business/articles/1809/13/news128.html",
    "https://www.itmedia.co.jp/business/articles/1808/27/news010.html"
)

$row = 3
foreach ($url in $remoteUrls) {
    $cell = $wsRemote.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row = $row + 1
}
